# Apply the changes described by the diff:
#  - Add a new student "Paulo Coelho Souza" (row 4), mirroring existing layout
#    but with blank/whitespace-only RA and TELEFONE values.
#  - Add a second new row (row 5) with a numeric RA/TELEFONE pair and
#    whitespace-only NOME/EMAIL values.
#  - Widen column B (NOME) slightly.
#  - Update the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new student with blank-looking RA/TELEFONE -----------------
# Order matches how the shared-string table ends up indexed in the
# target workbook (Paulo Coelho Souza, email, then the blank placeholders).
$ws.Range("B4").Value = "Paulo Coelho Souza"
$ws.Range("D4").Value = "paulo.souza@fatec.sp.gov.br"
$ws.Range("A4").Value = "     "
$ws.Range("C4").Value = "     "

# --- Row 5: numeric RA/TELEFONE pair with blank-looking NOME/EMAIL -----
$ws.Range("A5").Value = 1460311714078
$ws.Range("C5").Value = 1239665691
$ws.Range("B5").Value = "    "
$ws.Range("D5").Value = "    "

# Copy number formatting/styles from the row above (A3/C3) onto the new
# numeric cells (A5/C5) so they keep the same display style as the other
# RA / TELEFONE values in the table.
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial(-4122) # xlPasteFormats

# --- Column B width ------------------------------------------------------
# Widen the NOME column slightly so the new, longer name fits.
$ws.Columns.Item(2).ColumnWidth = 17.7

# --- Selection -------------------------------------------------------
$ws.Range("D5").Select()
